$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "Z00MN65"
$ws.Range("B2").Value = 44084.54782234977
$ws.Range("C2").Value = 44084.55539790984
$ws.Range("D2").Value = 10.9088065

# Delete row 3 entirely (shifts remaining rows up, reduces dimension)
$ws.Rows("3").Delete()
